# magius product type and product code
#
# Adds two new columns to the checklist sheet:
#   F: product_type (rulebook / supplement)
#   G: product_code (13-xx)
# and renames the sheet from the default "Sheet1" to "checklist".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "checklist"

# Header row
$ws.Range("F1").Value = "product_type"
$ws.Range("G1").Value = "product_code"

# row -> (product_type, product_code)   product_code "" means "no code for this row"
$data = @(
    @(2,  "rulebook",   "13-1"),
    @(3,  "supplement", "13-8"),
    @(4,  "supplement", "13-2"),
    @(5,  "supplement", "13-7"),
    @(6,  "supplement", "13-3"),
    @(7,  "supplement", "13-11"),
    @(8,  "supplement", "13-6"),
    @(9,  "supplement", "13-10"),
    @(10, "supplement", "13-9"),
    @(11, "supplement", "13-4"),
    @(12, "supplement", "13-5"),
    @(13, "supplement", "13-12"),
    @(14, "supplement", "13-15"),
    @(15, "supplement", "13-17"),
    @(16, "supplement", "13-13"),
    @(17, "supplement", "13-16"),
    @(18, "supplement", "13-14"),
    @(19, "supplement", ""),
    @(20, "supplement", "13-24"),
    @(21, "supplement", ""),
    @(22, "supplement", ""),
    @(23, "supplement", ""),
    @(24, "supplement", "13-23"),
    @(25, "supplement", ""),
    @(26, "supplement", ""),
    @(27, "supplement", "13-22"),
    @(28, "supplement", "13-28"),
    @(29, "supplement", ""),
    @(30, "supplement", "")
)

foreach ($entry in $data) {
    $row = $entry[0]
    $productType = $entry[1]
    $productCode = $entry[2]

    $ws.Range("F$row").Value = $productType
    if ($productCode -ne "") {
        $ws.Range("G$row").Value = $productCode
    }
}

# New column's width, matching the rest of the checklist's styling
$ws.Columns.Item(6).ColumnWidth = 14.832

# Leave the selection on the newly-added product_type column, mirroring the
# author's last selection after filling it in.
$ws.Range("F3:F30").Select()
